# Fix a vipLevel bug: refresh the "starttime" placeholder timestamp that was
# stamped when rows were generated, and correct a handful of rows whose
# TotalAward (M) / vipLevel (N) values were computed incorrectly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldStartTime = 1586862571.589257
$newStartTime = 1586918607.814806

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

# Column J is "starttime": every row still holding the stale placeholder
# timestamp gets refreshed to the new run's timestamp.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    if ($cell.Value() -eq $oldStartTime) {
        $cell.Value = $newStartTime
    }
}

# Column M is "TotalAward", column N is "vipLevel" - correct the mis-computed
# vipLevel (and its dependent TotalAward) for the affected accounts.
$ws.Cells.Item(39, 13).Value = 130.53
$ws.Cells.Item(39, 14).Value = 1

$ws.Cells.Item(93, 13).Value = 372.7

$ws.Cells.Item(147, 13).Value = 0
$ws.Cells.Item(147, 14).Value = 0

$ws.Cells.Item(224, 13).Value = 90
$ws.Cells.Item(224, 14).Value = 0

$ws.Cells.Item(228, 13).Value = 193.2795
